# Update countries & provincias Spain
# Applies the data refresh described by the diff:
#  - Irlanda overtakes Suecia in ranking (row 23/24 swap labels+values)
#  - Irak overtakes Islandia in ranking (row 69/70 swap labels+values)
#  - Several rows get refreshed statistic values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Francia) - refreshed stats
$ws.Range("D7").Value = 44903
$ws.Range("E7").Value = 93729
$ws.Range("F7").Value = 4682
$ws.Range("G7").Value = 242
$ws.Range("H7").Value = 22856

# Row 8 (Alemania) - refreshed stats
$ws.Range("D8").Value = 112000
$ws.Range("E8").Value = 39230

# Row 18 (Suiza) - refreshed stats
$ws.Range("E18").Value = 6151
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 1610

# Row 23/24: Irlanda overtakes Suecia
$ws.Range("A23").Value = "Irlanda"
$ws.Range("B23").Value = 19262
$ws.Range("C23").Value = 701
$ws.Range("D23").Value = 9233
$ws.Range("E23").Value = 8942
$ws.Range("F23").Value = 142
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 1087

$ws.Range("A24").Value = "Suecia"
$ws.Range("B24").Value = 18640
$ws.Range("C24").Value = 463
$ws.Range("D24").Value = 1005
$ws.Range("E24").Value = 15441
$ws.Range("F24").Value = 399
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 2194

# Row 69/70: Irak overtakes Islandia
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 1820
$ws.Range("C69").Value = 57
$ws.Range("D69").Value = 1263
$ws.Range("E69").Value = 470
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 87

$ws.Range("A70").Value = "Islandia"
$ws.Range("B70").Value = 1792
$ws.Range("C70").Value = 2
$ws.Range("D70").Value = 1608
$ws.Range("E70").Value = 174
$ws.Range("F70").Value = 3
$ws.Range("H70").Value = 10

# Row 91 (Bolivia) - refreshed stats
$ws.Range("D91").Value = 74
$ws.Range("E91").Value = 746

# Row 116 (Mali) - refreshed stats
$ws.Range("B116").Value = 389
$ws.Range("C116").Value = 19
$ws.Range("D116").Value = 112
$ws.Range("E116").Value = 254
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 23

# Row 170 (Siria) - refreshed stats
$ws.Range("B170").Value = 43
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 14
$ws.Range("E170").Value = 26
